# Aide.docx — "Modif aide et layout info"
#
# Applies, in bottom-to-top paragraph order (so earlier paragraph indices
# stay stable while later ones are touched first):
#
#   1. "S'inscrire" -> "A venir" gets a bold/red " (Est-ce necessaire?)" suffix.
#   2. New descriptive paragraph added right after "Profil".
#   3. New descriptive paragraph added right after "Statistiques".
#   4. New descriptive paragraph added right after "Mes dons".
#   5. "Faire un don" -> "A venir" body replaced with real descriptive text.
#   6. The _GoBack bookmark is relocated: removed from mid-sentence in the
#      "...crochet" paragraph and re-inserted as its own empty paragraph
#      right after the "Mes reservations" description.
#   7. "Les Organismes Communautaires" description: "par tous" -> "a tous".
#   8. "Connexion" -> "A venir" gets the same bold/red suffix as (1).

$d = $word.ActiveDocument

function Add-RedBoldSuffix($paragraphIndex) {
    $p = $d.Paragraphs($paragraphIndex)
    $r = $p.Range
    $r.InsertAfter(" ")
    $r2 = $p.Range
    $r2.InsertAfter("(Est-ce nécessaire?)")

    $full = $p.Range
    $suffix = "(Est-ce nécessaire?)"
    $start = $full.End - 1 - $suffix.Length
    $end = $full.End - 1
    $sub = $d.Range($start, $end)
    $sub.Font.Bold = 1
    $sub.Font.Color = 255
}

function Insert-BodyParagraphAfter($paragraphIndex, $text) {
    $p = $d.Paragraphs($paragraphIndex)
    $p.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs($paragraphIndex + 1)
    $newp.Style = "Normal"
    $newp.Range.Text = $text
}

# --- 1. "S'inscrire" section's "A venir" (last paragraph of that kind) ---
Add-RedBoldSuffix 32

# --- 2. "Profil" heading (paragraph 29) gains a description paragraph ---
Insert-BodyParagraphAfter 29 "C’est une fonctionnalité accessible à tous les utilisateurs inscrits."

# --- 3. "Statistiques" heading (paragraph 27) gains a description paragraph ---
Insert-BodyParagraphAfter 27 "C’est une fonctionnalité accessible à tous, excepté pour certaines informations."

# --- 4. "Mes dons" heading (paragraph 25) gains a description paragraph ---
Insert-BodyParagraphAfter 25 "C’est une fonctionnalité disponible à tous les utilisateurs inscrits et qui ont effectués des dons."

# --- 5. "Faire un don" body ("A venir", paragraph 24) becomes real text ---
$d.Paragraphs(24).Range.Text = "C’est une fonctionnalité disponible à tous les utilisateurs inscrits, puisque tout le monde est invité à faire des dons aux collectivités."

# --- 6. Relocate the _GoBack bookmark ---
# It currently sits mid-sentence, between the highlighted "crochet" run and
# the " » dans un cercle vert." run, inside the "La troisième affiche..."
# paragraph (paragraph 20). Remove it there...
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ...and recreate it as its own empty paragraph right after the
# "Mes réservations" description paragraph (paragraph 22).
$p22 = $d.Paragraphs(22)
$p22.Range.InsertParagraphAfter()
$bmPara = $d.Paragraphs(23)
$bmPara.Style = "Normal"
$bmRange = $bmPara.Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 7. "Les Organismes Communautaires" description: "par" -> "à" ---
$d.Content.Find.Execute("C’est une fonctionnalité accessible par tous.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C’est une fonctionnalité accessible à tous.", 2)

# --- 8. "Connexion" section's "A venir" (paragraph 7) ---
Add-RedBoldSuffix 7
